$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values updated to reflect refreshed crypto market data.
# Numeric-looking "Price" text values need NumberFormat="@" set first so
# Excel stores them verbatim as text (matching the source inlineStr cells)
# instead of auto-converting to floating point numbers.

$ws.Range('D2').Value = '22.398.53'
$ws.Range('E2').Value = '  -0.04%  '
$ws.Range('D3').Value = '1.572.46'
$ws.Range('E3').Value = '  +0.13%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '1.003'
$ws.Range('E5').Value = '  +0.18%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '290.23'
$ws.Range('E6').Value = '  -0.53%  '
$ws.Range('E7').Value = '  +3.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '49.93'
$ws.Range('E8').Value = '  +1.26%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3426'
$ws.Range('E9').Value = '  +0.92%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.159'
$ws.Range('E10').Value = '  -1.12%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07652'
$ws.Range('E11').Value = '  +0.93%  '
$ws.Range('E12').Value = '  +0.26%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.29'
$ws.Range('E13').Value = '  +0.51%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.022'
$ws.Range('E14').Value = '  -0.84%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.951'
$ws.Range('E15').Value = '  +0.83%  '
$ws.Range('D16').Value = '1.574.44'
$ws.Range('E16').Value = '  +0.14%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001134'
$ws.Range('E17').Value = '  -0.49%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '90.18'
$ws.Range('E18').Value = '  +1.19%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06762'
$ws.Range('E19').Value = '  +0.75%  '
$ws.Range('E20').Value = '  +0.13%  '
$ws.Range('E21').Value = '  +2.32%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.218'
$ws.Range('E22').Value = '  -0.57%  '
$ws.Range('E23').Value = '  -0.14%  '
$ws.Range('B24').Value = 'Toncoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.426'
$ws.Range('E24').Value = '  +1.55%  '
$ws.Range('B25').Value = 'WrappedBTC'
$ws.Range('C25').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D25').Value = '22.413.18'
$ws.Range('E25').Value = '  -0.01%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.720'
$ws.Range('E26').Value = '  -9.71%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.28'
$ws.Range('E27').Value = '  +1.90%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '147.43'
$ws.Range('E28').Value = '  +2.02%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.030'
$ws.Range('E29').Value = '  +1.06%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '126.42'
$ws.Range('E30').Value = '  +0.86%  '
$ws.Range('D31').Value = '1.750.83'
$ws.Range('E31').Value = '  -0.06%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.181'
$ws.Range('E32').Value = '  -1.59%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.009'
$ws.Range('E33').Value = '  +1.38%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.9910'
$ws.Range('E34').Value = '  -5.82%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '10.06'
$ws.Range('E35').Value = '  -3.38%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.08596'
$ws.Range('E36').Value = '  +0.99%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02550'
$ws.Range('E37').Value = '  -0.69%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2317'
$ws.Range('E38').Value = '  +0.61%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06577'
$ws.Range('E39').Value = '  +0.48%  '
$ws.Range('E40').Value = '  +5.78%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.483'
$ws.Range('E41').Value = '  -1.05%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6448'
$ws.Range('E42').Value = '  +0.89%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '11.54'
$ws.Range('E43').Value = '  -2.60%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.08'
$ws.Range('E44').Value = '  -3.87%  '
$ws.Range('B45').Value = 'Frax'
$ws.Range('C45').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.001'
$ws.Range('E45').Value = '  +0.17%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6011'
$ws.Range('E46').Value = '  -0.36%  '
$ws.Range('B47').Value = 'PancakeSwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.793'
$ws.Range('E47').Value = '  +0.26%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.307'
$ws.Range('E48').Value = '  +7.50%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.085'
$ws.Range('E49').Value = '  -2.33%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '125.62'
$ws.Range('E50').Value = '  +1.86%  '
$ws.Range('E51').Value = '  +0.40%  '
